# #5: property aircraft done
# Fix the property_category column (mistakenly "land" for every sheet) so
# that the 建物 (Building) sheet reports "building" and the 汽車 (Car) sheet
# reports "car" for their respective property rows.

$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet: column I is property_category, rows 2-5 hold data.
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I5").Value = "building"

# 汽車 (Car) sheet: column H is property_category, row 2 holds data.
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
